# feat: add 2022-Q4 data
#
# The workbook currently has 4 sheets: 总计, 2022-Q3, 2022-Q2, 2022-Q1.
# We need to insert a brand-new "2022-Q4" sheet (with the same layout as
# the existing quarter sheets) right before "2022-Q3", fill it with the
# new quarter's fund-holding data, and update the "总计" (totals) sheet
# so it lists all four quarters (Q4, Q3, Q2, Q1) with their counts and
# holding values.

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to be stored as text (matches the source workbook,
    # where these numeric-looking figures are inline/shared strings, not
    # real numbers), then drop back to the default "Normal" style so we
    # don't leave a stray number-format style on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (same
#    headers/column widths/styles), inserted right before it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with the Q4 fund-holding data.
# ---------------------------------------------------------------------
Set-TextValue $q4 "B2" "012868"
Set-TextValue $q4 "C2" "易方达标普信息科技指数（QDII-LOF）人民币 C"
Set-TextValue $q4 "D2" "5.09"
Set-TextValue $q4 "E2" "91.36"
Set-TextValue $q4 "F2" "1.86"
Set-TextValue $q4 "G2" "0.0947"
$q4.Range("H2").Value = 8

Set-TextValue $q4 "B3" "161128"
Set-TextValue $q4 "C3" "易方达标普信息科技指数（QDII-LOF）人民币"
Set-TextValue $q4 "D3" "5.09"
Set-TextValue $q4 "E3" "91.36"
Set-TextValue $q4 "F3" "1.86"
Set-TextValue $q4 "G3" "0.0947"
$q4.Range("H3").Value = 8

Set-TextValue $q4 "B4" "003721"
Set-TextValue $q4 "C4" "易方达标普信息科技指数（QDII-LOF）美元A"
Set-TextValue $q4 "D4" "4.93"
Set-TextValue $q4 "E4" "91.36"
Set-TextValue $q4 "F4" "1.86"
Set-TextValue $q4 "G4" "0.0917"
$q4.Range("H4").Value = 8

Set-TextValue $q4 "B5" "012869"
Set-TextValue $q4 "C5" "易方达标普信息科技指数（QDII-LOF）美元 C"
Set-TextValue $q4 "D5" "0.16"
Set-TextValue $q4 "E5" "91.36"
Set-TextValue $q4 "F5" "1.86"
Set-TextValue $q4 "G5" "0.0030"
$q4.Range("H5").Value = 8

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: it now has 4 rows (Q4, Q3, Q2, Q1)
#    instead of 3 (Q3, Q2, Q1). Extend the styled A-column down one row,
#    then rewrite the whole data block.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.28

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.28

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.2

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.23

# Keep "总计" as the active/selected sheet (unchanged from the source
# workbook), rather than leaving the newly-touched "2022-Q4" sheet active.
$total.Activate()
